# semana 22 de 2025
# Adds the new epidemiological week column "22" (column Y) to the weekly
# IRA-UCI revision sheet, mirroring the pattern of the existing week
# columns (1 .. 21 in D1:X1) and filling in the per-institution counts
# reported for that week (all zero except row 34, which reports 1; a
# handful of rows that have no data for any week are left blank, matching
# the source extract).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Y1 = "22" (stored as text, like the other week-number headers) ---
$ws.Range("Y1").Value = "'22"

# --- Body rows: per-institution count for week 22 ---
# Rows that already contain data for this institution this week (plain numbers).
$ws.Range("Y2").Value = 0

$ws.Range("Y5:Y8").Value = 0

$ws.Range("Y10:Y17").Value = 0

$ws.Range("Y19").Value = 0

$ws.Range("Y21:Y24").Value = 0

$ws.Range("Y27").Value = 0

$ws.Range("Y30:Y31").Value = 0

$ws.Range("Y33").Value = 0

$ws.Range("Y34").Value = 1

$ws.Range("Y36:Y49").Value = 0

$ws.Range("Y51:Y56").Value = 0

# Rows 3, 4, 9, 18, 20, 25, 26, 28, 29, 32, 35 and 50 have no entry for
# week 22 (same as in the source extract) and are intentionally left
# untouched.
